$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.07324766367673874
$ws.Range("B2").Value = 0.979211688041687
$ws.Range("C2").Value = 0.01092724315822124
$ws.Range("D2").Value = 0.9984169006347656
$ws.Range("A3").Value = 0.01145558338612318
$ws.Range("B3").Value = 0.9981393218040466
$ws.Range("C3").Value = 0.00556598138064146
$ws.Range("D3").Value = 0.9987184405326843
$ws.Range("A4").Value = 0.006082460749894381
$ws.Range("B4").Value = 0.9985243082046509
$ws.Range("C4").Value = 0.002481367439031601
$ws.Range("D4").Value = 0.9992461204528809
$ws.Range("A5").Value = 0.003369377693161368
$ws.Range("B5").Value = 0.9992728233337402
$ws.Range("C5").Value = 0.0007478162297047675
$ws.Range("D5").Value = 0.9998492002487183
$ws.Range("A6").Value = 0.003114670515060425
$ws.Range("B6").Value = 0.999315619468689
$ws.Range("C6").Value = 0.001410352764651179
$ws.Range("D6").Value = 0.9998492002487183
$ws.Range("A7").Value = 0.001158522441983223
$ws.Range("B7").Value = 0.9997219443321228
$ws.Range("C7").Value = 0.00005470982796396129
$ws.Range("D7").Value = 1
$ws.Range("A8").Value = 0.001284919213503599
$ws.Range("B8").Value = 0.9995936155319214
$ws.Range("C8").Value = 0.00007560154335806146
$ws.Range("D8").Value = 1
$ws.Range("A9").Value = 0.0005103262374177575
$ws.Range("B9").Value = 0.9998930692672729
$ws.Range("C9").Value = 0.00001219735440827208
$ws.Range("D9").Value = 1
$ws.Range("A10").Value = 0.001724665984511375
$ws.Range("B10").Value = 0.9996150135993958
$ws.Range("C10").Value = 0.00002891153781092726
$ws.Range("D10").Value = 1
$ws.Range("A11").Value = 0.001044031581841409
$ws.Range("B11").Value = 0.9997861385345459
$ws.Range("C11").Value = 0.001312674605287611
$ws.Range("D11").Value = 0.9998492002487183
$ws.Range("A12").Value = 0.0007520049693994224
$ws.Range("B12").Value = 0.9997006058692932
$ws.Range("C12").Value = 0.00002332657459191978
$ws.Range("D12").Value = 1
$ws.Range("A13").Value = 0.0005975269014015794
$ws.Range("B13").Value = 0.9998502731323242
$ws.Range("C13").Value = 0.000008368455382878892
$ws.Range("D13").Value = 1
$ws.Range("A14").Value = 0.0009217304759658873
$ws.Range("B14").Value = 0.9997219443321228
$ws.Range("C14").Value = 0.00004891554635833018
$ws.Range("D14").Value = 1
$ws.Range("A15").Value = 0.0004338714061304927
$ws.Range("B15").Value = 0.9998930692672729
$ws.Range("C15").Value = 0.000005277190211927518
$ws.Range("D15").Value = 1
$ws.Range("A16").Value = 0.0003039418661501259
$ws.Range("B16").Value = 0.9999358654022217
$ws.Range("C16").Value = 0.000001448734792575124
$ws.Range("D16").Value = 1
$ws.Range("A17").Value = 0.0007184931891970336
$ws.Range("B17").Value = 0.9997647404670715
$ws.Range("C17").Value = 0.00001811787842598278
$ws.Range("D17").Value = 1
$ws.Range("A18").Value = 0.0008103522704914212
$ws.Range("B18").Value = 0.9998502731323242
$ws.Range("C18").Value = 0.000002473867880325997
$ws.Range("D18").Value = 1
$ws.Range("A19").Value = 0.0007054029265418649
$ws.Range("B19").Value = 0.9998075366020203
$ws.Range("C19").Value = 0.000003653948851933819
$ws.Range("D19").Value = 1
$ws.Range("A20").Value = 0.0001667013420956209
$ws.Range("B20").Value = 0.9999358654022217
$ws.Range("C20").Value = 0.000000174530583763044
$ws.Range("D20").Value = 1
$ws.Range("A21").Value = 0.0003273676557000726
$ws.Range("B21").Value = 0.9999358654022217
$ws.Range("C21").Value = 0.00005579222488449886
$ws.Range("D21").Value = 1
$ws.Range("A22").Value = 0.0005812101298943162
$ws.Range("B22").Value = 0.9998716711997986
$ws.Range("C22").Value = 0.000006895519163663266
$ws.Range("D22").Value = 1
$ws.Range("A23").Value = 0.0002798312343657017
$ws.Range("B23").Value = 0.9999786019325256
$ws.Range("C23").Value = 0.0000001464009073970374
$ws.Range("D23").Value = 1
$ws.Range("A24").Value = 0.0003110703255515546
$ws.Range("B24").Value = 0.9999572038650513
$ws.Range("C24").Value = 0.00000007226390152936801
$ws.Range("D24").Value = 1
$ws.Range("A25").Value = 0.0002898909151554108
$ws.Range("B25").Value = 0.9999144673347473
$ws.Range("C25").Value = 0.0000002728988590661174
$ws.Range("D25").Value = 1
$ws.Range("A26").Value = 0.001501838560216129
$ws.Range("B26").Value = 0.9997861385345459
$ws.Range("C26").Value = 0.00000009229562181189976
$ws.Range("D26").Value = 1
$ws.Range("A27").Value = 0.0004498705384321511
$ws.Range("B27").Value = 0.9999144673347473
$ws.Range("C27").Value = 0.00000005370112887703726
$ws.Range("D27").Value = 1
$ws.Range("A28").Value = 0.0002218554873252288
$ws.Range("B28").Value = 0.9999144673347473
$ws.Range("C28").Value = 0.0000006460485906245594
$ws.Range("D28").Value = 1
$ws.Range("A29").Value = 0.000246458948822692
$ws.Range("B29").Value = 0.9999144673347473
$ws.Range("C29").Value = 0.00000002656297048986289
$ws.Range("D29").Value = 1
$ws.Range("A30").Value = 0.0001764651970006526
$ws.Range("B30").Value = 0.9999786019325256
$ws.Range("C30").Value = 0.000000008321621614015839
$ws.Range("D30").Value = 1
$ws.Range("A31").Value = 0.00002502496317902114
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 0.00000001851984876566348
$ws.Range("D31").Value = 1
$ws.Range("A32").Value = 0.0004969422589056194
$ws.Range("B32").Value = 0.9999144673347473
$ws.Range("C32").Value = 0.00000003800738923587232
$ws.Range("D32").Value = 1
$ws.Range("A33").Value = 0.0001198485551867634
$ws.Range("B33").Value = 0.9999786019325256
$ws.Range("C33").Value = 0.00000001756835743549345
$ws.Range("D33").Value = 1
$ws.Range("A34").Value = 0.0003767807211261243
$ws.Range("B34").Value = 0.9999144673347473
$ws.Range("C34").Value = 0.000001484997937950538
$ws.Range("D34").Value = 1
$ws.Range("A35").Value = 0.0003519279998727143
$ws.Range("B35").Value = 0.9999572038650513
$ws.Range("C35").Value = 0.00000000587728088419226
$ws.Range("D35").Value = 1
$ws.Range("A36").Value = 0.000008539104783267248
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 0.000000003432920836488051
$ws.Range("D36").Value = 1
$ws.Range("A37").Value = 0.0004567454161588103
$ws.Range("B37").Value = 0.9999144673347473
$ws.Range("C37").Value = 0.0000000002965626355599937
$ws.Range("D37").Value = 1
$ws.Range("A38").Value = 0.0001002631615847349
$ws.Range("B38").Value = 0.9999572038650513
$ws.Range("C38").Value = 0.000000006119931672543544
$ws.Range("D38").Value = 1
$ws.Range("A39").Value = 0.0001552142784930766
$ws.Range("B39").Value = 0.9999786019325256
$ws.Range("C39").Value = 0.0000001570708150211431
$ws.Range("D39").Value = 1
$ws.Range("A40").Value = 0.0008015321218408644
$ws.Range("B40").Value = 0.9998075366020203
$ws.Range("C40").Value = 0.0000000003325097974737901
$ws.Range("D40").Value = 1
$ws.Range("A41").Value = 0.0002073189534712583
$ws.Range("B41").Value = 0.9999572038650513
$ws.Range("C41").Value = 0.000000001662544768521457
$ws.Range("D41").Value = 1
$ws.Range("A42").Value = 0.0002846295246854424
$ws.Range("B42").Value = 0.9999358654022217
$ws.Range("C42").Value = 0.0000000003594700093145775
$ws.Range("D42").Value = 1
$ws.Range("A43").Value = 0.0005438351072371006
$ws.Range("B43").Value = 0.9999358654022217
$ws.Range("C43").Value = 0.0000000012042224994957
$ws.Range("D43").Value = 1
$ws.Range("A44").Value = 0.0008651061216369271
$ws.Range("B44").Value = 0.9999572038650513
$ws.Range("C44").Value = 0.0000000005931252711199875
$ws.Range("D44").Value = 1
$ws.Range("A45").Value = 0.0001517920318292454
$ws.Range("B45").Value = 0.9999572038650513
$ws.Range("C45").Value = 0.000000001051449038946828
$ws.Range("D45").Value = 1
$ws.Range("A46").Value = 0.000704568054061383
$ws.Range("B46").Value = 0.9998930692672729
$ws.Range("C46").Value = 0.000000004403489128890214
$ws.Range("D46").Value = 1
$ws.Range("A47").Value = 0.0002770397404674441
$ws.Range("B47").Value = 0.9999358654022217
$ws.Range("C47").Value = 0.000000004700050126871247
$ws.Range("D47").Value = 1
$ws.Range("A48").Value = 0.0004121703968849033
$ws.Range("B48").Value = 0.9998930692672729
$ws.Range("C48").Value = 0.000000002093908602063266
$ws.Range("D48").Value = 1
$ws.Range("A49").Value = 0.00006362393469316885
$ws.Range("B49").Value = 0.9999786019325256
$ws.Range("C49").Value = 0.0000000007099528187559656
$ws.Range("D49").Value = 1
$ws.Range("A50").Value = 0.00001471396899432875
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 0.0000000002156820555487826
$ws.Range("D50").Value = 1
$ws.Range("A51").Value = 0.000007940356226754375
$ws.Range("B51").Value = 1
$ws.Range("C51").Value = 0.0000000001527747789387135
$ws.Range("D51").Value = 1
